$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header cell - copy formatting (bold/border/alignment) from the
# neighboring header cell so it reuses the existing header style, then set
# its text.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Save column values (H2:H8)
$saveValues = @(0, 1, 0, 1, 0, 0, 0)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
